{"js": "// Commit: \"Delen van presentatie gemaakt\"\n// The run that starts with \"-Asset list + codes (assets + tijdsinschatting)\n// (paige)\" gets split: the bullet's own text is colored blue (004DBB) while\n// the rest of the paragraph (the following line breaks and remaining\n// bullets) keeps its original \"auto\" color.\n\nconst body = context.document.body;\n\nconst results = body.search(\n  \"-Asset list + codes (assets + tijdsinschatting) (paige)\",\n  { matchCase: true, matchWholeWord: false }\n);\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target text not found\");\n}\n\n// Color just the matched text blue; Word will split the owning run so the\n// rest of the paragraph (breaks + following bullets) remains untouched.\nresults.items[0].font.color = \"#004DBB\";\n\nawait context.sync();\n", "ps1": "# Commit: \"Delen van presentatie gemaakt\"\n# The run that starts with \"-Asset list + codes (assets + tijdsinschatting)\n# (paige)\" gets split: the bullet's own text is colored blue (RGB 00,4D,BB)\n# while the rest of the paragraph (the following line breaks and remaining\n# bullets) keeps its original \"auto\" color.\n\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$found = $rng.Find.Execute(\"-Asset list + codes (assets + tijdsinschatting) (paige)\")\n\nif ($found) {\n    # Word's Font.Color takes a BGR-packed long (0xBBGGRR) for RGB(0x00,0x4D,0xBB)\n    $rng.Font.Color = 0xBB4D00\n}\n"}
